$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the mobile-number columns are treated as text, not numbers
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M3").NumberFormat = "@"

# Update row 2 values
$ws.Range("L2").Value = "asdlkfjiozuxiojsdklfjj#!@dfax"
$ws.Range("M2").Value = "13918771256"
$ws.Range("N2").Value = "2022-12-29 00:00:00"
$ws.Range("Q2").Value = "system"

# Update row 3 values
$ws.Range("L3").Value = "asdlkfjiozuxiojsdklfjj#!fdsk32x"
$ws.Range("M3").Value = "13918771250"
$ws.Range("N3").Value = "2022-12-29 00:00:00"
$ws.Range("Q3").Value = "system"

# Remove column S entirely (shifts nothing else, column S had no header,
# only rows 2 and 3 carried data there) so the sheet dimension shrinks
# from A1:S3 to A1:R3
$ws.Range("S1:S3").EntireColumn.Delete()
